# Generate Report for handoff
# Updates the localization-status workbook to reflect a failed handoff:
#  - New source markdown filename (new GUID) replaces the old one everywhere
#  - Status text flips from "Ready for handoff" to "Handoff transform failed"
#  - Per-language sheets: the "Latest Handoff File" link is dropped (handoff
#    transform failed, so there is no handoff file), handoff datetime resets
#    to the zero value, and the handoff reason flips from "Include" to
#    "Ignored"

$wb = $excel.ActiveWorkbook

$newFile = "c6618b96-2854-42be-a23b-85bdba8859bd.md"
$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"
$newFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ece9fe606f202aa61ac2f887bfaa50736c2f4bad/e2e/c6618b96-2854-42be-a23b-85bdba8859bd.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ece9fe606f202aa61ac2f887bfaa50736c2f4bad/.localization-config"

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = $newStatus
$ws.Range("C2").Value = $newStatus

$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newFileUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config")

# ---- Per-language sheets ----------------------------------------------
$langs = @("zh-cn", "de-de")

foreach ($langName in $langs) {
    $ws = $wb.Worksheets.Item($langName)

    # Source file name (A2) + status (B2)
    $ws.Range("A2").Value = $newFile
    $ws.Range("B2").Value = $newStatus

    # Handoff failed -> no handoff file anymore, clear the whole cell
    $ws.Range("C2").Clear()

    # Handoff datetime resets to the zero value
    $ws.Range("D2").Value = $zeroDate

    # Handback datetime / reason
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("H2").Value = "Ignored"

    # Rebuild hyperlinks: drop the handoff-file link (C2 is gone), keep the
    # source file + config links pointing at the right display text
    $ws.Cells.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $newFileUrl, "", "", $newFile)
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config")
}
